# "Generate Report for Handoff"
#
# The localization status changes from "In Translation" to
# "Ready for handoff" and the associated timestamps are refreshed. This
# touches the three sheets that are kept in sync by the report generator:
#   - Overview   : zh-cn / de-de status cells + the "Latest HO Xliff
#                  Generate Date" (mirrors the de-de handoff timestamp)
#   - zh-cn      : Status + Latest Handoff Datetime
#   - de-de      : Status + Latest Handoff Datetime
#
# Because the Status text grows longer ("In Translation" -> "Ready for
# handoff"), the Status column is widened to fit the new text on every
# sheet that shows it.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-27 06:55:32"

# Widen the zh-cn / de-de status columns (E & F) to fit the new text.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-27 06:55:27"

# Widen the Status column (C) to fit the new text.
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-27 06:55:32"

# Widen the Status column (C) to fit the new text.
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
